# add create excel logmovement
# Update existing rows 2-21: RLDT (F) and OAORDT (X) move from 20250417 to
# 20250418, and OBSMCD (W) gains a "saleCode=" prefix. Then append a new
# data row 22 for the newly created order line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 21; $r++) {
    $ws.Range("F$r").Value = "'20250418"
    $ws.Range("X$r").Value = "'20250418"
    $ws.Range("W$r").Value = "saleCode=20277"
}

# New row 22 - empty text cells are written as a lone "'" (quote-prefix)
# so they land as empty-string Text cells rather than blank/null cells,
# matching the rest of the sheet's empty string cells.
$ws.Range("A22").Value = "'V01000390"
$ws.Range("B22").Value = "F10"
$ws.Range("C22").Value = "'215"
$ws.Range("D22").Value = "'"
$ws.Range("E22").Value = "'"
$ws.Range("F22").Value = "'20250418"
$ws.Range("G22").Value = "test"
$ws.Range("H22").Value = "'6804132150016"
$ws.Range("I22").Value = "'"
$ws.Range("J22").Value = "'10010101018"
$ws.Range("K22").Value = "'"
$ws.Range("L22").Value = "PCS"
$ws.Range("M22").Value = 3
$ws.Range("N22").Value = 109
$ws.Range("O22").Value = "PCS"
$ws.Range("P22").Value = "'"
$ws.Range("Q22").Value = "'"
$ws.Range("R22").Value = 1
$ws.Range("S22").Value = 0
$ws.Range("T22").Value = "'"
$ws.Range("U22").Value = "'"
$ws.Range("V22").Value = "'"
$ws.Range("W22").Value = "saleCode=20277"
$ws.Range("X22").Value = "'20250418"
$ws.Range("Y22").Value = "'"
$ws.Range("Z22").Value = "'"
$ws.Range("AA22").Value = "'"
$ws.Range("AB22").Value = "'"
$ws.Range("AC22").Value = "'"
$ws.Range("AD22").Value = "'"
$ws.Range("AE22").Value = "'"
$ws.Range("AF22").Value = "'"
$ws.Range("AG22").Value = "'"
$ws.Range("AH22").Value = "'"
$ws.Range("AI22").Value = "'"
$ws.Range("AJ22").Value = "'"
$ws.Range("AK22").Value = "'"
$ws.Range("AL22").Value = "'"
$ws.Range("AM22").Value = "'"
$ws.Range("AN22").Value = "'"
$ws.Range("AO22").Value = "'"
$ws.Range("AP22").Value = "'"
$ws.Range("AQ22").Value = "'"
$ws.Range("AR22").Value = "'"
$ws.Range("AS22").Value = "'"
